# Apply "primer abono aparece en historial abonos" changes.
$d = $word.ActiveDocument

# 1. The pre-existing "Fecha devolución" value (2019-02-27) becomes 2019-02-28.
#    Do this BEFORE the 2019-02-26 -> 2019-02-27 replace below, so we don't
#    accidentally touch the newly created 2019-02-27 values.
$d.Content.Find.Execute("2019-02-27", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "2019-02-28", 2)

# 2. Order/event dates: 2019-02-26 -> 2019-02-27 (both occurrences).
$d.Content.Find.Execute("2019-02-26", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "2019-02-27", 2)

# 3. Order number: No.  120 -> No.  124
$d.Content.Find.Execute("No.  120", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "No.  124", 2)

# 4. Abono history: first abono amount 10000 -> 29997
$d.Content.Find.Execute("10000", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "29997", 2)

# 5. Saldo amount: 40000 -> 20000
$d.Content.Find.Execute("40000", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "20000", 2)

# 6. TOTAL column value on the Saldo row: 50000 -> 49997.
#    "50000" also appears once in the item price cell above (which must stay
#    unchanged), so target the specific table cell instead of a global replace.
#    Use a Document.Range(...) (not the cell's own Range object) and
#    wdReplaceOne so the search/replace is confined to that single cell.
$t = $d.Tables.Item(2)
$cell = $t.Cell(7, 5)
$cellRange = $d.Range($cell.Range.Start, $cell.Range.End)
$cellRange.Find.Execute("50000", $false, $true, $false, $false, $false, `
                         $true, 0, $false, "49997", 1)
